# Edit script: split/move the AR "Unnamed: 43" status column into AS/AT/AU
# as per the commit "Used most updated status accomplishmnet files as of may"
#
# Summary of the change:
#  - Header row: AS1 becomes "Unnamed: 44", a new AT1 becomes "Unnamed: 45",
#    and the old AS1 header ("Status as of July 4, 2025") moves to AU1.
#  - Data rows (2-376): the old AR value (COMPLETED/ONGOING) moves to AT.
#    For some rows a new AS value ("ongrid") is inserted.
#    For other rows AR keeps a new "BBM ..." classification value while the
#    old AR value moves to AT.
#  - The dropdown data validation that used to target AS2:AS376 now targets
#    AU2:AU376.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows where a new AS value ("ongrid") must be inserted in addition to moving AR -> AT
$caseBRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,32,33,34,35,36,37,38,39,40,41,72,73,74,75)

# Rows where AR keeps a new value (old AR value moves to AT instead of being cleared)
$caseCValues = @{
    267 = "BBM 2025 UPGRADE"
    272 = "BBM 2025 UPGRADE"
    281 = "BBM 2025 UPGRADE"
    285 = "BBM 2025 UPGRADE"
    288 = "BBM 2025 UPGRADE"
    291 = "BBM 2025 UPGRADE"
    292 = "BBM 2025 UPGRADE"
    293 = "BBM 2025 UPGRADE"
    295 = "BBM 2023 UPGRADE"
    296 = "BBM 2025 UPGRADE"
    297 = "bbm 2023 ONGRID"
    311 = "BBM 2025 UPGRADE"
    318 = "BBM 2025 UPGRADE"
    326 = "bbm 2023 ONGRID"
    333 = "bbm 2023 ONGRID"
    338 = "BBM 2023 UPGRADE"
    339 = "bbm 2023 SOLAR"
    340 = "bbm 2023 SOLAR"
    341 = "BBM 2023 UPGRADE"
    342 = "BBM 2023 UPGRADE"
    343 = "bbm 2023 SOLAR"
    344 = "BBM 2023 UPGRADE"
    345 = "bbm 2023 SOLAR"
    346 = "BBM 2023 UPGRADE"
    347 = "bbm 2023 SOLAR"
    348 = "bbm 2023 SOLAR"
    349 = "BBM 2023 UPGRADE"
    350 = "BBM 2023 UPGRADE"
    351 = "BBM 2023 UPGRADE"
    352 = "BBM 2023 UPGRADE"
    353 = "BBM 2023 UPGRADE"
    354 = "BBM 2023 UPGRADE"
    355 = "BBM 2023 UPGRADE"
    356 = "BBM 2023 UPGRADE"
    357 = "BBM 2023 UPGRADE"
    358 = "BBM 2023 UPGRADE"
    359 = "BBM 2023 UPGRADE"
    360 = "BBM 2023 UPGRADE"
    361 = "BBM 2023 UPGRADE"
    362 = "bbm 2023 SOLAR"
    363 = "bbm 2023 SOLAR"
    366 = "BBM 2025 SOLAR"
    367 = "BBM 2024 UPGRADE"
    368 = "BBM 2024 SOLAR"
    369 = "BBM 2024 SOLAR"
    370 = "BBM 2024 SOLAR"
    371 = "BBM 2024 SOLAR"
    372 = "BBM 2024 SOLAR"
    373 = "BBM 2024 SOLAR"
    374 = "BBM 2025 ONGRID"
    375 = "BBM 2024 UPGRADE"
    376 = "BBM 2024 SOLAR"
}

# ---- Header row (row 1) ----
# AU1 gets the previous AS1 content/formatting (no special style)
$ws.Range("AU1").Value = $ws.Range("AS1").Value2

# AS1 and AT1 become new header labels with the bold/centered/bordered header style
$ws.Range("AS1").Value = "Unnamed: 44"
$ws.Range("AT1").Value = "Unnamed: 45"

$ws.Range("A1").Copy() | Out-Null
$ws.Range("AS1").PasteSpecial(-4122) | Out-Null
$ws.Range("AT1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---- Data rows (2-376) ----
for ($i = 2; $i -le 376; $i++) {
    $oldAR = $ws.Cells.Item($i, 44).Value2

    # Move old AR value to AT
    $ws.Cells.Item($i, 46).Value = $oldAR

    if ($caseBRows -contains $i) {
        $ws.Cells.Item($i, 45).Value = "ongrid"
        $ws.Cells.Item($i, 44).ClearContents() | Out-Null
    } elseif ($caseCValues.ContainsKey($i)) {
        $ws.Cells.Item($i, 44).Value = $caseCValues[$i]
    } else {
        $ws.Cells.Item($i, 44).ClearContents() | Out-Null
    }
}

# ---- Data validation: move from AS2:AS376 to AU2:AU376 ----
$ws.Range("AS2:AS376").Validation.Delete() | Out-Null
$ws.Range("AU2:AU376").Validation.Add(3, 1, 1, "=DropdownOptions!`$A`$1:`$A`$7") | Out-Null

Write-Output "done"
